# Update the "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column between "Week" (A) and "ASIN" (was B, now C)
#  - shorten the Week labels from "W01".."W16" to "W1".."W16"
#  - populate the new column with the week's start date (as text, Sundays
#    starting 2025-01-05, incrementing by 7 days)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column before column B (ASIN). This shifts every
# existing column from B..I to C..J (headers + all 16 data rows).
$ws.Columns.Item(2).Insert()

# New column header
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week start dates, one per data row (rows 2..17), Sundays 7 days apart
# starting 2025-01-05. Stored as literal text via a leading apostrophe so
# Excel doesn't reinterpret them as date serials.
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    # Shorten "W01".."W16" to "W1".."W16"
    $weekNum = $i + 1
    $ws.Cells.Item($row, 1).Value = "W" + $weekNum

    # Fill in the new Week_Start_Date column. Prefix with an apostrophe so
    # Excel stores the date as literal text instead of converting it to a
    # date serial number, then reset the style so no stray "quote prefix"
    # formatting is left behind on the cell.
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = "'" + $weekStartDates[$i]
    $cell.Style = "Normal"
}
